$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: paragraph 8 currently holds two runs ("4." and "Nakon ... pripada.")
# -- merge them into a single run with the combined text "4.Nakon ... pripada."
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$full = $p8.Range.Text
$full = $full.Substring(0, $full.Length - 1)   # drop the trailing paragraph mark
$splitIdx = 2                                   # length of "4."
$tailText = $full.Substring($splitIdx)

$run1Start = $p8.Range.Start
$run1End = $run1Start + $splitIdx
$run2Start = $run1End
$run2End = $p8.Range.End - 1

$tailRange = $d.Range($run2Start, $run2End)
$tailRange.Delete()

$headRange = $d.Range($run1Start, $run1End)
$headRange.InsertAfter($tailText)

# ---------------------------------------------------------------------------
# Step 2: helper to add a new "Normal" / justify-both paragraph right after a
# given paragraph, pre-populated with text (or empty).
# ---------------------------------------------------------------------------
function Add-ParaAfter($afterPara, [string]$text) {
    $afterPara.Range.InsertParagraphAfter()
    $newPara = $afterPara.Next()
    if ($text.Length -gt 0) {
        $newPara.Range.InsertBefore($text)
    }
    return $newPara
}

# ---------------------------------------------------------------------------
# Step 3: build the "Upiti:" section, paragraph by paragraph, right after the
# (now merged) paragraph 8.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(8)

$pUpiti = Add-ParaAfter $anchor "Upiti:"
$d.Bookmarks.Add("__DdeLink__14_370296924", $pUpiti.Range) | Out-Null

$p1c = Add-ParaAfter $pUpiti "1.c : U okviru ovog programa postoje 4 upita."

$pPrvi = Add-ParaAfter $p1c "-Prvi upit ispisuje id naziv I cenu stvari koje su prodate kupcu a kostale su vise od 2000 dinara."

$pDrugi = Add-ParaAfter $pPrvi "-Drugi upit brise stvar ciji je id=2."

$pTreci = Add-ParaAfter $pDrugi "-Treci upit unosi  novog dobavljaca u tabelu dobavljac."

$pCetvrti = Add-ParaAfter $pTreci "-Cetvrti upit pokusava da unese novu stvar u tabelu stvar I pri tome aktivira triger koji javlja da datum nije pravilno postavljen."

$pBlank1 = Add-ParaAfter $pCetvrti ""

$p2c = Add-ParaAfter $pBlank1 "2.c :"

$pOmogucava = Add-ParaAfter $p2c "-Omogucava dodavanje ili brisanje korisnika u zavisnosti koju opciju od ponudjenih izaberemo."

$p3c = Add-ParaAfter $pOmogucava "3.c:"

$pUnosimo = Add-ParaAfter $p3c "-Unosimo id stvari koja ce biti prodata kupcu I pri tome se dodaje novi red u tabelu prodata dok se u tabeli stvar kolicina te stvari umanjuje za 1.Nakon toga azurira se I stanje te stvari u jednoj od tabela obuca,odeca ili ostalo.(aktivira se triger3 I triger4)."
$d.Bookmarks.Add("__DdeLink__14_370296924", $pUnosimo.Range) | Out-Null

$pBlank2 = Add-ParaAfter $pUnosimo ""

Write-Output "done"
